# Update PSSM values in B2:K21 ("updates with supplemental figures").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,1.330826931665212,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,1.895158410817335,-15.88876702729936,3.851805963875106,-15.88876702729936,1.674704849149431,-15.88876702729936,2.553727013809188,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,3.565421800932354,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(3.019651982279626,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,1.734950173090626,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(3.571674870116053,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,0.4617751099304588,-15.88876702729936,1.965279527746481),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,2.168722379961399,-15.88876702729936,1.720736384597991,-15.88876702729936,-15.88876702729936,-15.88876702729936,1.300063251825594),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,2.711614155773321,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,1.140315256840086,2.233560506524214),
    @(-15.88876702729936,-15.88876702729936,1.333093820588632,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,2.538337078494662),
    @(-15.88876702729936,-15.88876702729936,0.8372515127626883,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,1.779543212543747,-15.88876702729936),
    @(-15.88876702729936,4.321905505415622,1.780374094569136,-15.88876702729936,-15.88876702729936,-15.88876702729936,1.775798310295907,1.92025369026608,2.526747816839961,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,-15.88876702729936,1.781757835440787,0.06844840070253129,1.44788771355895,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,2.630995904183428,-15.88876702729936,-15.88876702729936,-15.88876702729936,2.271188051775426,1.505971974293621,-15.88876702729936,-15.88876702729936),
    @(-15.88876702729936,-15.88876702729936,1.229210819063788,-15.88876702729936,2.475508722787403,-15.88876702729936,1.288428629532705,3.076118110591533,-15.88876702729936,1.64170789471106),
    @(-15.88876702729936,-15.88876702729936,-15.88876702729936,2.491664443568205,-15.88876702729936,2.282473309159232,1.419367544926118,-15.88876702729936,-15.88876702729936,-15.88876702729936),
)

for ($r = 0; $r -lt $values.Count; $r++) {
    $rowVals = $values[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
